# Update gh-pages output data (generated at a56beed)
# Applies refreshed "想去人数" (interest counts), one "是否有舞台" flag flip,
# and refreshed cover-image URLs across the 展览 and 全部类型 sheets,
# plus refreshed counts on the 演出 sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws4 = $wb.Worksheets.Item("全部类型")

# --- Sheet "展览" (sheet1) ---
$ws1.Range("F2").Value  = 1194
$ws1.Range("F4").Value  = 11220
$ws1.Range("F5").Value  = 1538
$ws1.Range("F6").Value  = 431
$ws1.Range("F7").Value  = 706
$ws1.Range("F8").Value  = 2198
$ws1.Range("F9").Value  = 648
$ws1.Range("F10").Value = 866
$ws1.Range("F11").Value = 468
$ws1.Range("H11").Value = $false
$ws1.Range("J11").Value = "//i2.hdslb.com/bfs/openplatform/202401/ghqTgq3H1705290829994.jpeg"
$ws1.Range("F12").Value = 347
$ws1.Range("F13").Value = 382
$ws1.Range("F14").Value = 351
$ws1.Range("F15").Value = 1119
$ws1.Range("F16").Value = 491
$ws1.Range("F17").Value = 873
$ws1.Range("F18").Value = 340
$ws1.Range("J18").Value = "//i0.hdslb.com/bfs/openplatform/202401/vGqHyZ1y1705290719000.jpeg"
$ws1.Range("F19").Value = 524
$ws1.Range("F20").Value = 824
$ws1.Range("F21").Value = 847
$ws1.Range("F22").Value = 66
$ws1.Range("F23").Value = 92
$ws1.Range("F24").Value = 90
$ws1.Range("F25").Value = 220
$ws1.Range("F26").Value = 558
$ws1.Range("F27").Value = 81
$ws1.Range("F28").Value = 35
$ws1.Range("F29").Value = 268

# --- Sheet "演出" (sheet2) ---
$ws2.Range("F3").Value = 685
$ws2.Range("F5").Value = 72
$ws2.Range("F7").Value = 763

# --- Sheet "全部类型" (sheet4) ---
$ws4.Range("F3").Value  = 1194
$ws4.Range("F4").Value  = 685
$ws4.Range("F7").Value  = 11220
$ws4.Range("F8").Value  = 1538
$ws4.Range("F9").Value  = 72
$ws4.Range("F10").Value = 431
$ws4.Range("F11").Value = 706
$ws4.Range("F12").Value = 2198
$ws4.Range("F13").Value = 648
$ws4.Range("F14").Value = 866
$ws4.Range("F16").Value = 468
$ws4.Range("H16").Value = $false
$ws4.Range("J16").Value = "//i2.hdslb.com/bfs/openplatform/202401/ghqTgq3H1705290829994.jpeg"
$ws4.Range("F17").Value = 347
$ws4.Range("F18").Value = 382
$ws4.Range("F19").Value = 351
$ws4.Range("F20").Value = 1119
$ws4.Range("F21").Value = 491
$ws4.Range("F22").Value = 763
$ws4.Range("F23").Value = 873
$ws4.Range("F24").Value = 340
$ws4.Range("J24").Value = "//i0.hdslb.com/bfs/openplatform/202401/vGqHyZ1y1705290719000.jpeg"
$ws4.Range("F25").Value = 524
$ws4.Range("F26").Value = 824
$ws4.Range("F27").Value = 847
$ws4.Range("F28").Value = 66
$ws4.Range("F29").Value = 92
$ws4.Range("F31").Value = 90
$ws4.Range("F32").Value = 220
$ws4.Range("F33").Value = 558
$ws4.Range("F34").Value = 81
$ws4.Range("F35").Value = 35
$ws4.Range("F36").Value = 268
